# TestData.xlsx touch-up:
#  - correct the "browser" value for the Amazon test row (Data!C6) from
#    "chromegrid" to "chrome"
#  - leave the workbook with the "Data" sheet active/selected at C6
#    (the "Test" sheet, where the cursor had been parked on D1, is no
#    longer the selected tab)

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("Test")
$wsData = $wb.Worksheets.Item("Data")

# Cursor starts out on "Test"!D1 (unchanged cell, just no longer the tab
# that ends up active/saved).
$wsTest.Activate()
$wsTest.Range("D1").Select()

# Move to "Data", fix the mis-typed browser name in C6, and leave the
# selection sitting on that cell.
$wsData.Activate()
$wsData.Range("C6").Select()
$wsData.Range("C6").Value = "chrome"
